$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F12").Value = 32568233.600000001
$ws.Range("F13").Value = 44875817.619999997
$ws.Range("F14").Value = 350000000
$ws.Range("F15").Value = -53616441.740000002
$ws.Range("F16").Value = 1451068426

$ws.Range("F18").Formula = "=SUM(F12:F17)"

$ws.Range("F19").Value = 1038368426

$ws.Range("F21").Formula = "=SUM(F18:F20)"

$ws.Range("F22").Value = 1026703455

$ws.Range("F26").Value = " "

$excel.Calculate()
